# Commit: "added exercies to slides 07"
#
# Two new "Ejercicio" (exercise) slides are inserted into the deck, using the
# same "Título y objetos" (Title and Content) layout already used by the
# existing "Ejercicio 07-01" slide (slide 4):
#
#   - "Ejercicio 07-02" is inserted right before the "Caso de estudio:" slide
#     (i.e. right after "Argumentos opcionales y por default").
#   - "Ejercicio 07-03" is inserted right after the "Caso de estudio:" slide
#     (i.e. right before the "IPython" title slide).
#
# Both new slides only get a title; their content placeholder is left blank,
# matching the pre-existing "Ejercicio 07-01" slide pattern.

$p = $ppt.ActivePresentation

# "Título y objetos" is the 2nd custom layout on the slide master - the same
# layout used by every regular content slide in this deck (e.g. slide 4,
# "Ejercicio 07-01").
$layout = $p.SlideMaster.CustomLayouts.Item(2)

# --- Insert "Ejercicio 07-03" first --------------------------------------
# Before this insertion, the deck still has its original 18 slides, with
# "Caso de estudio:" at position 9 and the "IPython" title slide at
# position 10. Inserting at position 10 places the new slide between them.
$s1 = $p.Slides.AddSlide(10, $layout)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Ejercicio 07-03"
$s1.Shapes.Item(1).TextFrame.TextRange.LanguageID = "es-MX"

# --- Insert "Ejercicio 07-02" second -------------------------------------
# The deck now has 19 slides; "Argumentos opcionales y por default" is still
# at position 8, and "Caso de estudio:" is still at position 9. Inserting at
# position 9 places the new slide right before "Caso de estudio:".
$s2 = $p.Slides.AddSlide(9, $layout)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Ejercicio 07-02"
$s2.Shapes.Item(1).TextFrame.TextRange.LanguageID = "es-MX"

Write-Output ("Slides after edit: " + $p.Slides.Count)
